$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.922.85'
$ws.Range("E2").Value = '  +1.62%  '
$ws.Range("D3").Value = '1.890.79'
$ws.Range("E3").Value = '  +1.52%  '
$ws.Range("E4").Value = '  -0.53%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.40%  '
$ws.Range("E7").Value = '  +0.63%  '
$ws.Range("E8").Value = '  +1.97%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07838'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9877'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  +1.93%  '
$ws.Range("D12").Value = '1.930.43'
$ws.Range("E12").Value = '  +3.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.052'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.25%  '
$ws.Range("E14").Value = '  +0.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06930'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009970'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.002'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.37%  '
$ws.Range("D21").Value = '28.938.94'
$ws.Range("E21").Value = '  +1.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.296'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.81%  '
$ws.Range("D24").Value = '2.175.18'
$ws.Range("E24").Value = '  +4.49%  '
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.05'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.83%  '
$ws.Range("E27").Value = '  +1.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.882'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.924'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.46'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9059'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.290'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.328'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.264'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.187'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05767'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.83%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02072'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.76%  '
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.735'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.50%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5675'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1768'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.734'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.290'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.81%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5354'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.07046'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("E48").Value = '  +2.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.527'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.065'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.16%  '
